# Updated cryptos list on Mon May 29 02:45:33 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column holds text (not real numbers - note the thousand-dot
# values like "28.308.81"), matching the source sheet's inlineStr cells.
# Prefixing with a leading apostrophe keeps Excel from reinterpreting
# numeric-looking text and silently dropping trailing zeros, e.g.
# "0.07820" -> 0.0782.
function Set-PriceText($row, $text) {
    $ws.Cells.Item($row, 4).Value = "'" + $text
}

function Set-VolumeText($row, $pct) {
    $ws.Cells.Item($row, 5).Value = "  " + $pct + "  "
}

# Simple per-row Price (D) / Volume-1h (E) refreshes -------------------------
$updates = @(
    @{Row=2;  D="28.308.81";    E="+3.12%"}
    @{Row=3;  D="1.923.44";     E="+2.70%"}
    @{Row=4;  D="1.008";        E="-0.81%"}
    @{Row=5;  D="318.24";       E="+1.58%"}
    @{Row=6;  D="1.007";        E="-0.79%"}
    @{Row=7;  D="0.4856";       E="+1.30%"}
    @{Row=8;  D="0.3853";       E="+2.35%"}
    @{Row=9;  D="0.07388";      E="+0.08%"}
    @{Row=10; D="0.9426";       E="+0.16%"}
    @{Row=11; D="20.95";        E="+1.06%"}
    @{Row=12; D="0.07820";      E="-0.94%"}
    @{Row=13; D="1.950.24";     E="+3.39%"}
    @{Row=14; D="5.531";        E="+1.68%"}
    @{Row=15; D="6.665";        E="+0.90%"}
    @{Row=16; D="91.63";        E="+0.83%"}
    @{Row=17; D="1.009";        E="-0.78%"}
    @{Row=18; D="0.000008878";  E="-0.23%"}
    @{Row=19; D="1.006";        E="-0.81%"}
    @{Row=20; D="28.335.92";    E="+3.04%"}
    @{Row=21; D="14.91";        E="+0.02%"}
    @{Row=22; D="5.171";        E="+0.47%"}
    @{Row=23; D="2.190.87";     E="+2.98%"}
    @{Row=24; D="10.97";        E="+2.47%"}
    @{Row=25; D="1.932";        E="-1.36%"}
    @{Row=26; D="156.38";       E="+1.28%"}
    @{Row=27; D="18.61";        E="+0.10%"}
    @{Row=28; D="2.103";        E="+4.40%"}
    @{Row=29; D="116.67";       E="+0.47%"}
    @{Row=30; D="4.988";        E="-0.49%"}
    @{Row=31; D="0.08930";      E="-0.03%"}
    @{Row=32; D="3.370";        E="+1.17%"}
    @{Row=33; D="1.249";        E="+2.79%"}
    @{Row=34; D="0.7750";       E="+3.66%"}
    @{Row=35; D="4.711";        E="+2.94%"}
    @{Row=36; D="2.715";        E="+0.87%"}
    @{Row=37; D="0.02055";      E="-0.55%"}
    @{Row=38; D="1.112";        E="-1.11%"}
    @{Row=39; D="0.5584";       E="+4.25%"}
    @{Row=40; D="0.05342";      E="+0.63%"}
    @{Row=41; D="3.025";        E="+0.83%"}
    @{Row=42; D="7.071";        E=$null}
    @{Row=47; D="107.12";       E="+3.66%"}
    @{Row=48; D="1.008";        E="-0.85%"}
    @{Row=49; D="1.665";        E="+0.08%"}
    @{Row=50; D="69.28";        E="+3.11%"}
    @{Row=51; D="0.06147";      E="+0.43%"}
)

foreach ($u in $updates) {
    Set-PriceText $u.Row $u.D
    if ($u.E -ne $null) {
        Set-VolumeText $u.Row $u.E
    }
}

# Rows 43-46 were re-ranked: the coin name/link/price/volume move to a new
# row while the row index (col A) and styling stay put. Net effect is a
# swap of (43<->44) and (45<->46) content blocks with updated figures.
$reranked = @(
    @{Row=43; B="Algorand";     C="https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo";     D="0.1532"; E="+0.17%"}
    @{Row=44; B="Aptos";        C="https://coinranking.com/coin/HGYj5JCv5+aptos-apt";              D="8.474";  E="+0.58%"}
    @{Row=45; B="EnergySwap";   C="https://coinranking.com/coin/SbWqqTui-+energyswap-ens";          D="10.76";  E="+1.01%"}
    @{Row=46; B="Decentraland"; C="https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana";   D="0.4886"; E="+0.84%"}
)

foreach ($u in $reranked) {
    $ws.Cells.Item($u.Row, 2).Value = $u.B
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    Set-PriceText $u.Row $u.D
    Set-VolumeText $u.Row $u.E
}
